# Generate Report for Handoff
# Adds a new localization entry (file 8a1b7149-7071-4066-8814-3de76465f0f2)
# as row 3 on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$guid = "8a1b7149-7071-4066-8814-3de76465f0f2"
$hash = "7051c9a158f5f74d62c702eb7d967f02488d84bc"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = "$guid.md"
$wsOverview.Cells.Item(3, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 4).Value = "2016-27-14 03:27:35"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/8725b1402a4143afdeeddc3440056647a8b80dff/e2e/$guid.md",
    "",
    "",
    "$guid.md"
)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(3, 1).Value = "$guid.md"
$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = "$guid.$hash.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 5).Value = "2016-03-14 03:27:33"
$wsZhCn.Cells.Item(3, 8).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3, 9).Value = "Include"

$wsZhCn.Cells.Item(3, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/8725b1402a4143afdeeddc3440056647a8b80dff/e2e/$guid.md",
    "",
    "",
    "$guid.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(3, 2),
    "https://github.com/OpenLocalizationTest/oltest/blob/8725b1402a4143afdeeddc3440056647a8b80dff/e2e/$guid.md",
    "",
    "",
    ".md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(3, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$guid.$hash.zh-cn.xlf",
    "",
    "",
    "$guid.$hash.zh-cn.xlf"
)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(3, 1).Value = "$guid.md"
$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = "$guid.$hash.de-de.xlf"
$wsDeDe.Cells.Item(3, 5).Value = "2016-03-14 03:27:35"
$wsDeDe.Cells.Item(3, 8).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3, 9).Value = "Include"

$wsDeDe.Cells.Item(3, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(3, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/8725b1402a4143afdeeddc3440056647a8b80dff/e2e/$guid.md",
    "",
    "",
    "$guid.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(3, 2),
    "https://github.com/OpenLocalizationTest/oltest/blob/8725b1402a4143afdeeddc3440056647a8b80dff/e2e/$guid.md",
    "",
    "",
    ".md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(3, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$hash/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$guid.$hash.de-de.xlf",
    "",
    "",
    "$guid.$hash.de-de.xlf"
)
